$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027335366613918
$ws.Range("D2").Value = 1.034881058108952
$ws.Range("E2").Value = 1.027387393856686
$ws.Range("F2").Value = 1.042010425207702
$ws.Range("I2").Value = 1.023594999628091
$ws.Range("J2").Value = 1.032493811895247
$ws.Range("K2").Value = 1.037679260548263
$ws.Range("L2").Value = 1.030207270465784
$ws.Range("M2").Value = 1.044788321311973
$ws.Range("N2").Value = 1.033960071536852

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028727479668758
$ws.Range("D3").Value = 1.036174129765367
$ws.Range("E3").Value = 1.028585289557306
$ws.Range("F3").Value = 1.043473244770953
$ws.Range("I3").Value = 1.023504579208683
$ws.Range("J3").Value = 1.033523793130703
$ws.Range("K3").Value = 1.038780324607272
$ws.Range("L3").Value = 1.031211825966324
$ws.Range("M3").Value = 1.046060176221304
$ws.Range("N3").Value = 1.034991515463785

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029627962063276
$ws.Range("D4").Value = 1.037010819290603
$ws.Range("E4").Value = 1.029360361741198
$ws.Range("F4").Value = 1.044420050316333
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.034189557042108
$ws.Range("K4").Value = 1.039492230414934
$ws.Range("L4").Value = 1.031861234186112
$ws.Range("M4").Value = 1.046882906415014
$ws.Range("N4").Value = 1.035658224836308

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030006456032151
$ws.Range("D5").Value = 1.037362563931431
$ws.Range("E5").Value = 1.029686194208471
$ws.Range("F5").Value = 1.044818154975199
$ws.Range("I5").Value = 1.023416932628374
$ws.Range("J5").Value = 1.034469279179316
$ws.Range("K5").Value = 1.039791385866541
$ws.Range("L5").Value = 1.032134103209044
$ws.Range("M5").Value = 1.047228727531681
$ws.Range("N5").Value = 1.035938344211035

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030070002940092
$ws.Range("D6").Value = 1.03742162353571
$ws.Range("E6").Value = 1.029740902545407
$ws.Range("F6").Value = 1.044885002595281
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.034516236174856
$ws.Range("K6").Value = 1.03984160787391
$ws.Range("L6").Value = 1.032179910866667
$ws.Range("M6").Value = 1.047286789329637
$ws.Range("N6").Value = 1.035985367890895

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029633019782616
$ws.Range("D7").Value = 1.037015519313573
$ws.Range("E7").Value = 1.029364715557514
$ws.Range("F7").Value = 1.044425369538983
$ws.Range("I7").Value = 1.023442985510381
$ws.Range("J7").Value = 1.034193295347354
$ws.Range("K7").Value = 1.039496228250306
$ws.Range("L7").Value = 1.03186488083039
$ws.Range("M7").Value = 1.046887527503508
$ws.Range("N7").Value = 1.035661968450377

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027805904410675
$ws.Range("D8").Value = 1.035318062972578
$ws.Range("E8").Value = 1.027792239882794
$ws.Range("F8").Value = 1.042504740103736
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.032842046354945
$ws.Range("K8").Value = 1.038051487541045
$ws.Range("L8").Value = 1.03054689204197
$ws.Range("M8").Value = 1.045218203762036
$ws.Range("N8").Value = 1.034308800529407

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024583718032479
$ws.Range("D9").Value = 1.032326648677793
$ws.Range("E9").Value = 1.025020835944708
$ws.Range("F9").Value = 1.039122160523541
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.03045544471044
$ws.Range("K9").Value = 1.035501256147781
$ws.Range("L9").Value = 1.028219650109888
$ws.Range("M9").Value = 1.042274586937265
$ws.Range("N9").Value = 1.031918809636821

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022433556368672
$ws.Range("D10").Value = 1.030331945821979
$ws.Range("E10").Value = 1.023172691857872
$ws.Range("F10").Value = 1.036868042639149
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.028860455335135
$ws.Range("K10").Value = 1.03379791867725
$ws.Range("L10").Value = 1.026664755764723
$ws.Range("M10").Value = 1.040310542663705
$ws.Range("N10").Value = 1.030321555193655

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021501961062663
$ws.Range("D11").Value = 1.029468062673378
$ws.Range("E11").Value = 1.022372250587149
$ws.Range("F11").Value = 1.035892142708855
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.028168831217058
$ws.Range("K11").Value = 1.033059551325002
$ws.Range("L11").Value = 1.025990620762614
$ws.Range("M11").Value = 1.039459646155185
$ws.Range("N11").Value = 1.02962894889

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021155834842576
$ws.Range("D12").Value = 1.02914714873944
$ws.Range("E12").Value = 1.022074899577408
$ws.Range("F12").Value = 1.035529666700692
$ws.Range("I12").Value = 1.023938164268129
$ws.Range("J12").Value = 1.027911779487192
$ws.Range("K12").Value = 1.032785162737396
$ws.Range("L12").Value = 1.025740084775817
$ws.Range("M12").Value = 1.039143512426875
$ws.Range("N12").Value = 1.029371532117193

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021230084271147
$ws.Range("D13").Value = 1.029215987267473
$ws.Range("E13").Value = 1.022138683868013
$ws.Range("F13").Value = 1.035607418369639
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.027966924887085
$ws.Range("K13").Value = 1.032844025794328
$ws.Range("L13").Value = 1.025793831638828
$ws.Range("M13").Value = 1.039211327561784
$ws.Range("N13").Value = 1.029426755829879

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021473352036747
$ws.Range("D14").Value = 1.029441536451703
$ws.Range("E14").Value = 1.022347672125561
$ws.Range("F14").Value = 1.035862180020938
$ws.Range("I14").Value = 1.023922712353276
$ws.Range("J14").Value = 1.028147586359624
$ws.Range("K14").Value = 1.033036872875254
$ws.Range("L14").Value = 1.025969914098746
$ws.Range("M14").Value = 1.03943351595779
$ws.Range("N14").Value = 1.029607673862431

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021623225199578
$ws.Range("D15").Value = 1.029580500729455
$ws.Range("E15").Value = 1.022476432395781
$ws.Range("F15").Value = 1.036019149125891
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.028258877628021
$ws.Range("K15").Value = 1.033155675493107
$ws.Range("L15").Value = 1.026078386680774
$ws.Range("M15").Value = 1.039570403757194
$ws.Range("N15").Value = 1.029719123177194

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022495370952211
$ws.Range("D16").Value = 1.03038927505438
$ws.Range("E16").Value = 1.023225810315807
$ws.Range("F16").Value = 1.03693281254774
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.028906335092203
$ws.Range("K16").Value = 1.033846904238986
$ws.Range("L16").Value = 1.026709477544684
$ws.Range("M16").Value = 1.04036700391255
$ws.Range("N16").Value = 1.03036750010524

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023042290525656
$ws.Range("D17").Value = 1.030896550708069
$ws.Range("E17").Value = 1.023695823532942
$ws.Range("F17").Value = 1.0375059651788
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.029312201989157
$ws.Range("K17").Value = 1.034280273659996
$ws.Range("L17").Value = 1.027105112521694
$ws.Range("M17").Value = 1.040866565566367
$ws.Range("N17").Value = 1.030773943379742

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023361245723389
$ws.Range("D18").Value = 1.031192420542773
$ws.Range("E18").Value = 1.02396995693842
$ws.Range("F18").Value = 1.037840290067522
$ws.Range("I18").Value = 1.023826118601784
$ws.Range("J18").Value = 1.029548842692929
$ws.Range("K18").Value = 1.034532972750474
$ws.Range("L18").Value = 1.027335797419519
$ws.Range("M18").Value = 1.041157908113013
$ws.Range("N18").Value = 1.031010920140456

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023469992267763
$ws.Range("D19").Value = 1.031293302027153
$ws.Range("E19").Value = 1.024063426523684
$ws.Range("F19").Value = 1.037954288868121
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.029629515134557
$ws.Range("K19").Value = 1.034619123514042
$ws.Range("L19").Value = 1.027414441161802
$ws.Range("M19").Value = 1.04125724119365
$ws.Range("N19").Value = 1.031091707146204

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022983616834215
$ws.Range("D20").Value = 1.030842126466818
$ws.Range("E20").Value = 1.023645397388189
$ws.Range("F20").Value = 1.037444469851583
$ws.Range("I20").Value = 1.023846096030192
$ws.Range("J20").Value = 1.029268666126463
$ws.Range("K20").Value = 1.034233785318249
$ws.Range("L20").Value = 1.02706267317257
$ws.Range("M20").Value = 1.040812971845881
$ws.Range("N20").Value = 1.030730345691129

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021401718282093
$ws.Range("D21").Value = 1.029375118684259
$ws.Range("E21").Value = 1.022286131194175
$ws.Range("F21").Value = 1.03578715864874
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.028094390279943
$ws.Range("K21").Value = 1.032980087730962
$ws.Range("L21").Value = 1.025918065883094
$ws.Range("M21").Value = 1.039368089064449
$ws.Range("N21").Value = 1.029554402238215

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020406590792212
$ws.Range("D22").Value = 1.028452580996316
$ws.Range("E22").Value = 1.021431322430792
$ws.Range("F22").Value = 1.034745235042274
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.027355196534637
$ws.Range("K22").Value = 1.032191106744529
$ws.Range("L22").Value = 1.025197638956192
$ws.Range("M22").Value = 1.038459209809445
$ws.Range("N22").Value = 1.028814158753004

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020934178464495
$ws.Range("D23").Value = 1.028941653519684
$ws.Range("E23").Value = 1.021884491537471
$ws.Range("F23").Value = 1.035297571435951
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.027747141878148
$ws.Range("K23").Value = 1.032609431254527
$ws.Range("L23").Value = 1.025579624834041
$ws.Range("M23").Value = 1.038941065744512
$ws.Range("N23").Value = 1.02920666070386

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02301012911545
$ws.Range("D24").Value = 1.030866718486321
$ws.Range("E24").Value = 1.023668182845848
$ws.Range("F24").Value = 1.037472256897304
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.02928833840087
$ws.Range("K24").Value = 1.034254791639458
$ws.Range("L24").Value = 1.027081849941721
$ws.Range("M24").Value = 1.040837188674797
$ws.Range("N24").Value = 1.030750045902422

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025417067088078
$ws.Range("D25").Value = 1.033100058618121
$ws.Range("E25").Value = 1.025737390675232
$ws.Range("F25").Value = 1.039996452642754
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.031073114602026
$ws.Range("K25").Value = 1.036161097295538
$ws.Range("L25").Value = 1.028219650109888
$ws.Range("M25").Value = 1.04303585180819
$ws.Range("N25").Value = 1.032537356690501

Write-Host "Applied vm_pu updates for rows 2-25 (380 kV case)"
